$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Change B11 value from "R40" to the text "1" (kept as text, not the number 1),
# matching the new shared-string entry added by the source edit. A leading
# apostrophe forces Excel to store this as text rather than a numeric value.
$ws.Range("B11").Value = "'1"
